# Portugal Primeira Liga - odds base update (2024-01-29 10:56)
# 1) Nine pairs of adjacent rows had their match-record contents (everything
#    except the running id in column A) swapped with one another.
# 2) Five still-upcoming fixtures had their odds columns (N,O,P,Q,R,S,U,V)
#    refreshed with newer quotes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

$rowPairs = @(
    @(667, 668),
    @(695, 696),
    @(722, 723),
    @(742, 743),
    @(752, 753),
    @(776, 777),
    @(778, 779),
    @(866, 867),
    @(982, 983)
)

foreach ($pair in $rowPairs) {
    Swap-Rows $pair[0] $pair[1]
}

# Direct odds refresh for upcoming fixtures (no result yet recorded)
$ws.Range("N987").Value2 = 10
$ws.Range("O987").Value2 = 6.5
$ws.Range("P987").Value2 = 1.25
$ws.Range("Q987").Value2 = 1.75
$ws.Range("R987").Value2 = 1.86
$ws.Range("S987").Value2 = 2.04
$ws.Range("U987").Value2 = 1.8
$ws.Range("V987").Value2 = 2.05

$ws.Range("N988").Value2 = 1.181
$ws.Range("O988").Value2 = 7.5
$ws.Range("P988").Value2 = 13
$ws.Range("Q988").Value2 = -2
$ws.Range("R988").Value2 = 1.88
$ws.Range("S988").Value2 = 2.02
$ws.Range("U988").Value2 = 2.05
$ws.Range("V988").Value2 = 1.8

$ws.Range("R989").Value2 = 1.93
$ws.Range("S989").Value2 = 1.97
$ws.Range("U989").Value2 = 1.95
$ws.Range("V989").Value2 = 1.9

$ws.Range("N991").Value2 = 3
$ws.Range("P991").Value2 = 2.3
$ws.Range("Q991").Value2 = 0.25
$ws.Range("R991").Value2 = 1.87
$ws.Range("S991").Value2 = 2.03
$ws.Range("U991").Value2 = 1.825
$ws.Range("V991").Value2 = 2.025

$ws.Range("N994").Value2 = 3.3
$ws.Range("O994").Value2 = 3.4
$ws.Range("P994").Value2 = 2.15
$ws.Range("R994").Value2 = 2.03
$ws.Range("S994").Value2 = 1.87
$ws.Range("U994").Value2 = 1.875
$ws.Range("V994").Value2 = 1.975
